$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oktober")

# Fill in the hours and details for the newly logged days (rows 24-26)
# Order matches the order new shared strings were appended in the file.
$ws.Range("B24").Value = 2.5
$ws.Range("C24").Value = "CC2650 in betrieb genommen, TI RTOS eingelesen, "

$ws.Range("B26").Value = 1.5
$ws.Range("C26").Value = "CC2650 I2C Beispiel laufen lassen, nun spinnt der debugger "

$ws.Range("B25").Value = 2
$ws.Range("C25").Value = "I2C Beispiel geschrieben, gibt fehler beim start der Transaction "

# Update the active selection to reflect where the user ended up
$ws.Activate()
$ws.Range("C30").Select()
